$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) was shared
# between the 6361649d-... (row 3) and a5f462f1-... (row 4) rows, so both
# get updated together.
$wsOverview.Range("G3").Value = "2016-08-25 14:14:36"
$wsOverview.Range("G4").Value = "2016-08-25 14:14:36"

# zh-cn sheet: rows 3 (6361649d-...) and 4 (a5f462f1-...) shared the same
# Priority / Correspond Handoff Datetime / Correspond Handback DateTime values.
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-25 14:14:32"
$wsZhCn.Range("K3").Value = "2016-08-25 14:14:49"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-25 14:14:32"
$wsZhCn.Range("K4").Value = "2016-08-25 14:14:49"

# de-de sheet: rows 3 (6361649d-...) and 4 (a5f462f1-...) shared the same
# Priority / Correspond Handoff Datetime / Correspond Handback DateTime values.
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-25 14:14:36"
$wsDeDe.Range("K3").Value = "2016-08-25 14:14:55"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-25 14:14:36"
$wsDeDe.Range("K4").Value = "2016-08-25 14:14:55"
